$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 0.12249740311861
$ws.Range("C16").Value = 1.68190507979382
$ws.Range("D16").Value = 0.1365196219199632
$ws.Range("E16").Value = 1.403362708735537
